$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 91; this shifts all existing rows 91..127 down to 92..128
$ws.Rows.Item(91).Insert()

# Populate the newly inserted row 91 with the new record's data
$ws.Cells.Item(91, 1).Value = 8
$ws.Cells.Item(91, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(91, 3).Value = "Coquimbo"
$ws.Cells.Item(91, 4).Value = 44722
$ws.Cells.Item(91, 5).Value = 4
$ws.Cells.Item(91, 6).Value = "Fruta"
$ws.Cells.Item(91, 7).Value = 100109
$ws.Cells.Item(91, 8).Value = "Uva"
$ws.Cells.Item(91, 9).Value = 100109001
$ws.Cells.Item(91, 10).Value = "Uva"
$ws.Cells.Item(91, 11).Value = "Red Globe"
$ws.Cells.Item(91, 12).Value = "Primera"
$ws.Cells.Item(91, 13).Value = 248
$ws.Cells.Item(91, 14).Value = 9500
$ws.Cells.Item(91, 15).Value = 10000
$ws.Cells.Item(91, 16).Value = 9758
$ws.Cells.Item(91, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(91, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(91, 19).Value = 542
$ws.Cells.Item(91, 20).Value = 18

# Ensure the D column of the new row keeps the date-like numeric format used by the rest of column D
$ws.Cells.Item(91, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
